$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.934.81'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '3.543.29'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.45'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.69'
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("D7").Value = '3.538.66'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.612'
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.198'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.24'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.75'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").Value = '4.113.67'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '612.05'
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").Value = '3.538.03'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = '70.994.90'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.81'
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.888'
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.05'
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.74'
$ws.Range("E24").Value = '  -2.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.56'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.79'
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.92'
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("E30").Value = '  +0.99%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.88'
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '632.06'
$ws.Range("E35").Value = '  +9.71%  '
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.51'
$ws.Range("E38").Value = '  -5.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0477'
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.00'
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.144'
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("D43").Value = '0.0₃0741'
$ws.Range("E43").Value = '  +5.07%  '
$ws.Range("D44").Value = '3.368.18'
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.01'
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.24'
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.57'
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.15'
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.157'
$ws.Range("E51").Value = '  +6.44%  '
